$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.53
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 2.2
$ws.Range("L2").Value = 8
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38
$ws.Range("W2").Value = 4.5
$ws.Range("X2").Value = 5.5
$ws.Range("Z2").Value = 10
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 8
$ws.Range("AL2").Value = 81
$ws.Range("AN2").Value = 3.2
$ws.Range("AO2").Value = 8
$ws.Range("AT2").Value = 2.38
$ws.Range("BA2").Value = 301
$ws.Range("G3").Value = 4.33
$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 2.1
$ws.Range("O3").Value = 1.67
$ws.Range("P3").Value = 2.1
$ws.Range("S3").Value = 1.73
$ws.Range("T3").Value = 2.08
$ws.Range("AC3").Value = 5
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 8
$ws.Range("AJ3").Value = 19
$ws.Range("AK3").Value = 23
$ws.Range("AX3").Value = 13
$ws.Range("BB3").Value = 401
$ws.Range("BD5").Value = 151
$ws.Range("J6").Value = 1.83
$ws.Range("K6").Value = 2.6
$ws.Range("L6").Value = 6.5
$ws.Range("S6").Value = 1.25
$ws.Range("T6").Value = 3.75
$ws.Range("AO6").Value = 6.5
$ws.Range("AP6").Value = 15
$ws.Range("AQ6").Value = 17
$ws.Range("AR6").Value = 34
$ws.Range("AS6").Value = 101
$ws.Range("AT6").Value = 3.75
$ws.Range("AW6").Value = 9
$ws.Range("AX6").Value = 34
$ws.Range("AZ6").Value = 126
$ws.Range("BA6").Value = 126
